# Se agrego columna fecha emitida en reporte de compra
# The F2 cell (Año Modelo Vehículo Motor) held the stray text
# "CHEVROLET CAVALIER 2019" instead of just the model year, unlike every
# other row in the table which stores a plain numeric year. Fix it to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 2019

# Reflect the cell the user ended up on after making the edit.
$ws.Range("F3").Select()
